$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 145, shifting existing row 145 (and below) down to 146
$ws.Rows.Item(145).Insert()

# Populate the newly inserted row 145 with the new record's values
$ws.Cells.Item(145, 1).Value = 3
$ws.Cells.Item(145, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(145, 3).Value = "Coquimbo"
$ws.Cells.Item(145, 4).Value = 44900
$ws.Cells.Item(145, 4).NumberFormat = $ws.Cells.Item(146, 4).NumberFormat
$ws.Cells.Item(145, 5).Value = 5
$ws.Cells.Item(145, 6).Value = 100112052
$ws.Cells.Item(145, 7).Value = "Albahaca"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 120
$ws.Cells.Item(145, 11).Value = 6500
$ws.Cells.Item(145, 12).Value = 7000
$ws.Cells.Item(145, 13).Value = 6792
$ws.Cells.Item(145, 14).Value = "$/docena de matas"
$ws.Cells.Item(145, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(145, 16).Value = 1132
$ws.Cells.Item(145, 17).Value = 6
$ws.Cells.Item(145, 18).Value = "Hortaliza"
